$d = $word.ActiveDocument

# Update the date heading in the first paragraph
$d.Content.Find.Execute("2023-05-22 Monday", $true, $true, $false, $false, $false, $true, 1, $false, "2023-05-23 Tuesday", 2) | Out-Null

# Update each arithmetic-problem cell in the 20x5 table, by position,
# so the substring relationships between expressions (e.g. "7+41=" is
# contained in "37+41=") can never cause a wrong/double replacement.
$t = $d.Tables.Item(1)
$values = @(
    "77-37=",
    "71-51=",
    "92-40=",
    "15+46=",
    "15+65=",
    "70-7=",
    "44+20=",
    "58+18=",
    "11+74=",
    "82-22=",
    "42-21=",
    "23+62=",
    "55-24=",
    "41-30=",
    "66+13=",
    "74-16=",
    "37-21=",
    "85-55=",
    "13+61=",
    "36-0=",
    "41+11=",
    "91-23=",
    "54-8=",
    "25+39=",
    "76-11=",
    "12+12=",
    "69+1=",
    "38+4=",
    "12+67=",
    "53+35=",
    "44-14=",
    "35+16=",
    "83-33=",
    "70-47=",
    "22+57=",
    "97-4=",
    "61-5=",
    "28-20=",
    "69-0=",
    "32+51=",
    "45-25=",
    "15+41=",
    "36+48=",
    "99-88=",
    "1+13=",
    "0+0=",
    "62-46=",
    "6+70=",
    "71-31=",
    "67-17=",
    "0+41=",
    "41+9=",
    "89-45=",
    "98-85=",
    "26-8=",
    "34+36=",
    "94-3=",
    "20+33=",
    "34+3=",
    "28-5=",
    "13+70=",
    "16+54=",
    "41-13=",
    "27+48=",
    "95-76=",
    "32+13=",
    "99-86=",
    "10-3=",
    "59-22=",
    "31+38=",
    "92-80=",
    "37+19=",
    "5+93=",
    "66-38=",
    "74+13=",
    "99-52=",
    "77-15=",
    "55-28=",
    "44-9=",
    "95-4=",
    "89-6=",
    "30+62=",
    "7+89=",
    "44-24=",
    "57-1=",
    "11-0=",
    "82-63=",
    "24+21=",
    "57-7=",
    "13+9=",
    "88-33=",
    "21-20=",
    "75-40=",
    "35-17=",
    "23-10=",
    "59+34=",
    "40+57=",
    "87-6=",
    "82-35=",
    "95-21="
)

$rows = 20
$cols = 5
$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $values[$idx]
        $idx++
    }
}
